$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Segments")
$ws.Activate()

# Flip the boolean "Active" flag (column O) from FALSE to TRUE for rows 2-7
$ws.Range("O2:O7").Value = $true

# Update the active cell / selection to K4 and scroll the view back to show column A
$ws.Range("K4").Select()
